# Auto-generated edit script: update crypto price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.456.91'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +1.09%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.476.44'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -1.27%  '

# Row 4
$ws.Range("E4").Value = '  +0.38%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '488.15'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.24%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.73'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +6.76%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.15%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.510'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.07%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.488.26'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.63%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.78'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +6.18%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0978'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.79%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.334'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +2.19%  '

# Row 13
$ws.Range("E13").Value = '  +1.93%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.920.51'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.21%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '56.461.80'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.13%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.01'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.84%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000136'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.92%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.500.78'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.16%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.53'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +4.60%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.17'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +2.12%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '318.10'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.91%  '

# Row 22
$ws.Range("E22").Value = '  +0.29%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.87'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +3.62%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '58.47'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +1.15%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.409'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.64%  '

# Row 26
$ws.Range("B26").Value = 'Binance-PegBSC-USD'
$ws.Range("C26").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.80%  '

# Row 27
$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.163'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.51%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.599.51'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.78%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.59'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +2.60%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0788'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +2.80%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.19%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '148.84'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.66%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.26'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.84%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.51'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +2.06%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.19'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.12%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.16'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +5.78%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.75'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +1.71%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.872'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +3.32%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.38'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +4.07%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '33.84'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -1.87%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.49'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +3.59%  '

# Row 42
$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.996'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.22%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0555'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.85%  '

# Row 44
$ws.Range("B44").Value = 'Mantle'
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.609'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.87%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.85'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +10.60%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '264.54'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +5.34%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0925'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.78%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0229'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +3.19%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.23'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.99%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.69'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.43%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.890.64'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -4.08%  '
